$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 143.44444
$ws.Cells.Item(9, 9).Value = 65.333336
$ws.Cells.Item(9, 10).Value = 299.66666
$ws.Cells.Item(9, 11).Value = 65.333336
$ws.Cells.Item(9, 12).Value = 299.66666
$ws.Cells.Item(9, 13).Value = 103.666664
$ws.Cells.Item(9, 14).Value = -637.66666
$ws.Cells.Item(19, 8).Value = 952.3570999999999
$ws.Cells.Item(19, 9).Value = 684.8
$ws.Cells.Item(19, 10).Value = 1621.25
$ws.Cells.Item(19, 11).Value = 684.8
$ws.Cells.Item(19, 12).Value = 1621.25
$ws.Cells.Item(19, 13).Value = -509.8
$ws.Cells.Item(19, 14).Value = -1971.25
$ws.Cells.Item(38, 8).Value = 708.8333
$ws.Cells.Item(87, 8).Value = 15935.118
$ws.Cells.Item(87, 10).Value = 15935.118
$ws.Cells.Item(87, 12).Value = 15935.118
$ws.Cells.Item(87, 14).Value = -18431.118
$ws.Cells.Item(90, 8).Value = 15935.118
$ws.Cells.Item(90, 10).Value = 15935.118
$ws.Cells.Item(90, 12).Value = 47805.354
$ws.Cells.Item(90, 14).Value = -60285.354
$ws.Cells.Item(121, 8).Value = 1173.2084
$ws.Cells.Item(121, 10).Value = 1248.1
$ws.Cells.Item(121, 12).Value = 3744.3
$ws.Cells.Item(121, 14).Value = -7238.299999999999
$ws.Cells.Item(129, 8).Value = 924.0833
$ws.Cells.Item(129, 10).Value = 929.0571
$ws.Cells.Item(129, 12).Value = 2787.1713
$ws.Cells.Item(129, 14).Value = -12787.1713
$ws.Cells.Item(132, 8).Value = 14310.514
$ws.Cells.Item(132, 9).Value = 16182.375
$ws.Cells.Item(132, 10).Value = 2330.6
$ws.Cells.Item(132, 11).Value = 48547.125
$ws.Cells.Item(132, 12).Value = 6991.799999999999
$ws.Cells.Item(132, 13).Value = -46017.125
$ws.Cells.Item(132, 14).Value = -12051.8
$ws.Cells.Item(135, 8).Value = 2581657.5
$ws.Cells.Item(135, 9).Value = 3323.3845
$ws.Cells.Item(135, 10).Value = 36100000
$ws.Cells.Item(135, 11).Value = 29910.4605
$ws.Cells.Item(135, 12).Value = 324900000
$ws.Cells.Item(135, 13).Value = -27375.4605
$ws.Cells.Item(135, 14).Value = -324905070
$ws.Cells.Item(137, 8).Value = 1162.1471
$ws.Cells.Item(137, 9).Value = 1087.5652
$ws.Cells.Item(137, 10).Value = 1318.091
$ws.Cells.Item(137, 11).Value = 3262.6956
$ws.Cells.Item(137, 12).Value = 3954.273
$ws.Cells.Item(137, 13).Value = -712.6956
$ws.Cells.Item(137, 14).Value = -9054.272999999999
$ws.Cells.Item(138, 8).Value = 3966.4658
$ws.Cells.Item(138, 9).Value = 2457.7632
$ws.Cells.Item(138, 10).Value = 5113.08
$ws.Cells.Item(138, 11).Value = 7373.2896
$ws.Cells.Item(138, 12).Value = 15339.24
$ws.Cells.Item(138, 13).Value = -2233.2896
$ws.Cells.Item(138, 14).Value = -25619.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1977.8572
$ws.Cells.Item(61, 9).Value = 1732.4
$ws.Cells.Item(61, 10).Value = 2591.5
$ws.Cells.Item(61, 11).Value = 1732.4
$ws.Cells.Item(61, 12).Value = 2591.5
$ws.Cells.Item(61, 13).Value = -1520.4
$ws.Cells.Item(61, 14).Value = -3015.5
$ws.Cells.Item(74, 8).Value = 1542.0819
$ws.Cells.Item(74, 9).Value = 977.7234
$ws.Cells.Item(74, 10).Value = 3436.7144
$ws.Cells.Item(74, 11).Value = 977.7234
$ws.Cells.Item(74, 12).Value = 3436.7144
$ws.Cells.Item(74, 13).Value = -103.7234
$ws.Cells.Item(74, 14).Value = -5184.7144
$ws.Cells.Item(77, 8).Value = 1542.0819
$ws.Cells.Item(77, 9).Value = 977.7234
$ws.Cells.Item(77, 10).Value = 3436.7144
$ws.Cells.Item(77, 11).Value = 4888.617
$ws.Cells.Item(77, 12).Value = 17183.572
$ws.Cells.Item(77, 13).Value = -520.6170000000002
$ws.Cells.Item(77, 14).Value = -25919.572
$ws.Cells.Item(132, 8).Value = 2122.1333
$ws.Cells.Item(132, 9).Value = 1674.3334
$ws.Cells.Item(132, 10).Value = 3167
$ws.Cells.Item(132, 11).Value = 5023.0002
$ws.Cells.Item(132, 12).Value = 9501
$ws.Cells.Item(132, 13).Value = -2493.0002
$ws.Cells.Item(132, 14).Value = -14561
$ws.Cells.Item(134, 8).Value = 43610
$ws.Cells.Item(134, 10).Value = 48332
$ws.Cells.Item(134, 12).Value = 48332
$ws.Cells.Item(134, 14).Value = -58472
$ws.Cells.Item(136, 8).Value = 1977.8572
$ws.Cells.Item(136, 9).Value = 1732.4
$ws.Cells.Item(136, 10).Value = 2591.5
$ws.Cells.Item(136, 11).Value = 5197.200000000001
$ws.Cells.Item(136, 12).Value = 7774.5
$ws.Cells.Item(136, 13).Value = -2647.200000000001
$ws.Cells.Item(136, 14).Value = -12874.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 302
$ws.Cells.Item(22, 9).Value = 218.8
$ws.Cells.Item(22, 10).Value = 510
$ws.Cells.Item(22, 11).Value = 218.8
$ws.Cells.Item(22, 12).Value = 510
$ws.Cells.Item(22, 13).Value = -45.80000000000001
$ws.Cells.Item(22, 14).Value = -856
$ws.Cells.Item(133, 8).Value = 49898.9
$ws.Cells.Item(133, 10).Value = 49898.9
$ws.Cells.Item(133, 12).Value = 49898.9
$ws.Cells.Item(133, 14).Value = -60018.9
$ws.Cells.Item(140, 8).Value = 58668
$ws.Cells.Item(140, 10).Value = 58668
$ws.Cells.Item(140, 12).Value = 58668
$ws.Cells.Item(140, 14).Value = -69028

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4373.5933
$ws.Cells.Item(31, 9).Value = 2361.3167
$ws.Cells.Item(31, 10).Value = 9017.308000000001
$ws.Cells.Item(31, 11).Value = 2361.3167
$ws.Cells.Item(31, 12).Value = 9017.308000000001
$ws.Cells.Item(31, 13).Value = -2066.3167
$ws.Cells.Item(31, 14).Value = -9607.308000000001
$ws.Cells.Item(34, 8).Value = 4373.5933
$ws.Cells.Item(34, 9).Value = 2361.3167
$ws.Cells.Item(34, 10).Value = 9017.308000000001
$ws.Cells.Item(34, 11).Value = 2361.3167
$ws.Cells.Item(34, 12).Value = 9017.308000000001
$ws.Cells.Item(34, 13).Value = -2159.3167
$ws.Cells.Item(34, 14).Value = -9421.308000000001
$ws.Cells.Item(99, 8).Value = 1530892.1
$ws.Cells.Item(99, 9).Value = 2668983
$ws.Cells.Item(99, 10).Value = 13437.556
$ws.Cells.Item(99, 11).Value = 2668983
$ws.Cells.Item(99, 12).Value = 13437.556
$ws.Cells.Item(99, 13).Value = -2667485
$ws.Cells.Item(99, 14).Value = -16433.556
$ws.Cells.Item(122, 8).Value = 1835821
$ws.Cells.Item(122, 9).Value = 334137.34
$ws.Cells.Item(122, 11).Value = 1002412.02
$ws.Cells.Item(122, 13).Value = -999962.02
$ws.Cells.Item(126, 8).Value = 1530892.1
$ws.Cells.Item(126, 9).Value = 2668983
$ws.Cells.Item(126, 10).Value = 13437.556
$ws.Cells.Item(126, 11).Value = 8006949
$ws.Cells.Item(126, 12).Value = 40312.66800000001
$ws.Cells.Item(126, 13).Value = -8004479
$ws.Cells.Item(126, 14).Value = -45252.66800000001
$ws.Cells.Item(132, 8).Value = 677591.4
$ws.Cells.Item(132, 9).Value = 962590.9
$ws.Cells.Item(132, 11).Value = 2887772.7
$ws.Cells.Item(132, 13).Value = -2885242.7
$ws.Cells.Item(140, 8).Value = 74231.78
$ws.Cells.Item(140, 10).Value = 74231.78
$ws.Cells.Item(140, 12).Value = 74231.78
$ws.Cells.Item(140, 14).Value = -84591.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1600.3077
$ws.Cells.Item(113, 9).Value = 2337.8333
$ws.Cells.Item(113, 10).Value = 968.1429000000001
$ws.Cells.Item(113, 11).Value = 7013.499899999999
$ws.Cells.Item(113, 12).Value = 2904.4287
$ws.Cells.Item(113, 13).Value = -4843.499899999999
$ws.Cells.Item(113, 14).Value = -7244.4287
$ws.Cells.Item(122, 8).Value = 938.6829
$ws.Cells.Item(122, 9).Value = 556.6786
$ws.Cells.Item(122, 10).Value = 1761.4615
$ws.Cells.Item(122, 11).Value = 5010.1074
$ws.Cells.Item(122, 12).Value = 15853.1535
$ws.Cells.Item(122, 13).Value = -2560.1074
$ws.Cells.Item(122, 14).Value = -20753.1535
$ws.Cells.Item(131, 8).Value = 779.16
$ws.Cells.Item(131, 10).Value = 843
$ws.Cells.Item(131, 12).Value = 2529
$ws.Cells.Item(131, 14).Value = -12609

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 5131809.5
$ws.Cells.Item(102, 9).Value = 8550251
$ws.Cells.Item(102, 10).Value = 4146.6665
$ws.Cells.Item(102, 11).Value = 8550251
$ws.Cells.Item(102, 12).Value = 4146.6665
$ws.Cells.Item(102, 13).Value = -8548629
$ws.Cells.Item(102, 14).Value = -7390.6665
$ws.Cells.Item(122, 8).Value = 64600.9
$ws.Cells.Item(122, 9).Value = 82192.91
$ws.Cells.Item(122, 10).Value = 4006.2222
$ws.Cells.Item(122, 11).Value = 246578.73
$ws.Cells.Item(122, 12).Value = 12018.6666
$ws.Cells.Item(122, 13).Value = -244128.73
$ws.Cells.Item(122, 14).Value = -16918.6666
$ws.Cells.Item(126, 8).Value = 2239.2258
$ws.Cells.Item(126, 9).Value = 2126.3333
$ws.Cells.Item(126, 10).Value = 2395.5386
$ws.Cells.Item(126, 11).Value = 6378.999899999999
$ws.Cells.Item(126, 12).Value = 7186.6158
$ws.Cells.Item(126, 13).Value = -3908.999899999999
$ws.Cells.Item(126, 14).Value = -12126.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 146840.86
$ws.Cells.Item(7, 9).Value = 203577.6
$ws.Cells.Item(7, 10).Value = 4999
$ws.Cells.Item(7, 11).Value = 203577.6
$ws.Cells.Item(7, 12).Value = 4999
$ws.Cells.Item(7, 13).Value = -203465.6
$ws.Cells.Item(7, 14).Value = -5223
$ws.Cells.Item(122, 8).Value = 2391.5
$ws.Cells.Item(122, 9).Value = 2333.8
$ws.Cells.Item(122, 11).Value = 7001.400000000001
$ws.Cells.Item(122, 13).Value = -4551.400000000001
$ws.Cells.Item(126, 8).Value = 146840.86
$ws.Cells.Item(126, 9).Value = 203577.6
$ws.Cells.Item(126, 10).Value = 4999
$ws.Cells.Item(126, 11).Value = 610732.8
$ws.Cells.Item(126, 12).Value = 14997
$ws.Cells.Item(126, 13).Value = -608262.8
$ws.Cells.Item(126, 14).Value = -19937
$ws.Cells.Item(127, 8).Value = 60250
$ws.Cells.Item(127, 10).Value = 60250
$ws.Cells.Item(127, 12).Value = 60250
$ws.Cells.Item(127, 14).Value = -70170
$ws.Cells.Item(132, 8).Value = 17558
$ws.Cells.Item(132, 9).Value = 6217.3335
$ws.Cells.Item(132, 10).Value = 28898.666
$ws.Cells.Item(132, 11).Value = 18652.0005
$ws.Cells.Item(132, 12).Value = 86695.99800000001
$ws.Cells.Item(132, 13).Value = -16122.0005
$ws.Cells.Item(132, 14).Value = -91755.99800000001
$ws.Cells.Item(133, 8).Value = 41663
$ws.Cells.Item(133, 10).Value = 41663
$ws.Cells.Item(133, 12).Value = 41663
$ws.Cells.Item(133, 14).Value = -46723

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(128, 8).Value = 49630
$ws.Cells.Item(128, 10).Value = 49630
$ws.Cells.Item(128, 12).Value = 49630
$ws.Cells.Item(128, 14).Value = -59590
$ws.Cells.Item(136, 8).Value = 5873.9614
$ws.Cells.Item(136, 9).Value = 4186.9165
$ws.Cells.Item(136, 10).Value = 7320
$ws.Cells.Item(136, 11).Value = 12560.7495
$ws.Cells.Item(136, 12).Value = 21960
$ws.Cells.Item(136, 13).Value = -10010.7495
$ws.Cells.Item(136, 14).Value = -27060
$ws.Cells.Item(137, 8).Value = 55074.4
$ws.Cells.Item(137, 10).Value = 55074.4
$ws.Cells.Item(137, 12).Value = 55074.4
$ws.Cells.Item(137, 14).Value = -65274.4
